# Update cryptocurrency price/volume data per Sun Oct 27 2024 GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.925.93'
$ws.Range('E2').Value = '  +1.07%  '
$ws.Range('D3').Value = '2.504.36'
$ws.Range('E3').Value = '  +0.57%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '588.91'
$ws.Range('E5').Value = '  +0.76%  '
$ws.Range('D6').Value = '177.59'
$ws.Range('E6').Value = '  +2.82%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '0.515'
$ws.Range('E8').Value = '  +0.50%  '
$ws.Range('E9').Value = '  +3.26%  '
$ws.Range('E11').Value = '  +2.32%  '
$ws.Range('D12').Value = '4.95'
$ws.Range('E12').Value = '  +0.61%  '
$ws.Range('D13').Value = '2.957.04'
$ws.Range('E13').Value = '  +0.54%  '
$ws.Range('D14').Value = '25.71'
$ws.Range('E14').Value = '  +1.13%  '
$ws.Range('D15').Value = '67.700.58'
$ws.Range('E15').Value = '  +0.73%  '
$ws.Range('E16').Value = '  +0.78%  '
$ws.Range('D17').Value = '2.521.14'
$ws.Range('E17').Value = '  +1.56%  '
$ws.Range('D18').Value = '11.00'
$ws.Range('E18').Value = '  -0.12%  '
$ws.Range('D19').Value = '7.51'
$ws.Range('E19').Value = '  +1.18%  '
$ws.Range('D20').Value = '352.73'
$ws.Range('E20').Value = '  +0.96%  '
$ws.Range('E21').Value = '  +1.90%  '
$ws.Range('E22').Value = '  +0.14%  '
$ws.Range('D23').Value = '70.92'
$ws.Range('E23').Value = '  +3.44%  '
$ws.Range('D24').Value = '4.30'
$ws.Range('E24').Value = '  +1.57%  '
$ws.Range('D25').Value = '1.75'
$ws.Range('E25').Value = '  -1.83%  '
$ws.Range('D26').Value = '9.14'
$ws.Range('E26').Value = '  -1.32%  '
$ws.Range('D27').Value = '2.590.28'
$ws.Range('E27').Value = '  -0.95%  '
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  -0.01%  '
$ws.Range('D29').Value = '0.0₃0913'
$ws.Range('E29').Value = '  +1.17%  '
$ws.Range('D30').Value = '505.26'
$ws.Range('E30').Value = '  -0.78%  '
$ws.Range('D31').Value = '7.86'
$ws.Range('E31').Value = '  +1.10%  '
$ws.Range('E32').Value = '  +2.95%  '
$ws.Range('D33').Value = '1.78'
$ws.Range('E33').Value = '  +0.70%  '
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').Value = '0.122'
$ws.Range('E35').Value = '  +3.67%  '
$ws.Range('B36').Value = 'Monero'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D36').Value = '164.45'
$ws.Range('E36').Value = '  +2.92%  '
$ws.Range('D37').Value = '18.66'
$ws.Range('E37').Value = '  -0.25%  '
$ws.Range('D38').Value = '18.40'
$ws.Range('E38').Value = '  +0.92%  '
$ws.Range('E39').Value = '  +0.19%  '
$ws.Range('E40').Value = '  +0.06%  '
$ws.Range('E41').Value = '  +2.82%  '
$ws.Range('E42').Value = '  +0.56%  '
$ws.Range('D43').Value = '4.89'
$ws.Range('E43').Value = '  +1.29%  '
$ws.Range('E44').Value = '  +4.54%  '
$ws.Range('D45').Value = '145.18'
$ws.Range('E45').Value = '  +1.71%  '
$ws.Range('E46').Value = '  +2.80%  '
$ws.Range('D47').Value = '0.519'
$ws.Range('E47').Value = '  +0.89%  '
$ws.Range('D48').Value = '0.0₆0257'
$ws.Range('E48').Value = '  +2.69%  '
$ws.Range('D49').Value = '0.0744'
$ws.Range('E49').Value = '  +1.70%  '
$ws.Range('E50').Value = '  +1.37%  '
$ws.Range('D51').Value = '0.588'
$ws.Range('E51').Value = '  +0.82%  '
